$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.255.93'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.92%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.795.22'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.92%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '326.91'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.51%  '

$ws.Range("E6").Value = '  +0.08%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4519'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +15.20%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3748'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +10.25%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '44.73'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.35%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.146'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.28%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07543'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.14%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.58'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.13%  '

$ws.Range("E13").Value = '  +0.02%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.296'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.24%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.543'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +6.14%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.789.84'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.73%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001090'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.92%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06721'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.53%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '81.09'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.68%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.0000'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.01%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.56'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.57%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.350'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.96%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.240.54'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.80%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.77'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.93%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.421'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.25%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '20.54'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.76%  '

$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '151.56'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.33%  '

$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.352'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.44%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.977.64'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.90%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.09'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.87%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.232'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.43%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.017'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.52%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09428'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +8.18%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.819'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.15%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.2324'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +10.05%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.12'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.24%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06340'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.67%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02328'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.75%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.169'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.57%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6558'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.75%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.311'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.63%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.471'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.83%  '

$ws.Range("E43").Value = '  +0.27%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9994'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.01%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.05'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.39%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6099'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.51%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.791'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.07%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '129.94'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.42%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.026'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.58%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07128'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.80%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.160'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.19%  '
